# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-22) was re-sorted so the
# periods run in ascending order (2404 -> 2410, oldest to newest) instead of
# descending (2410 -> 2404). Periods 2405-2409 (rows 17-21) already sit in
# the right order, so only the first and last rows of the range need to
# trade places: row 16 (was period 2410 / 15600) <-> row 22 (was period
# 2404 / 52000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodoTop = $ws.Range("E16").Value()
$valorTop   = $ws.Range("F16").Value()

$periodoBottom = $ws.Range("E22").Value()
$valorBottom   = $ws.Range("F22").Value()

$ws.Range("E16").Value = $periodoBottom
$ws.Range("F16").Value = $valorBottom

$ws.Range("E22").Value = $periodoTop
$ws.Range("F22").Value = $valorTop
